# fix network overview LL layer.
#
# Slide 2 ("Implementation / Network architecture-and-modules_minimize")
# had its "Graphical Engine" abstraction-layer band (a filled rounded
# rectangle + its underline connector) and the "C stacks TCP/IP, TLS,
# Crypto, Bluetooth" band removed, and replaced with a single, wider
# band reading "C stacks TCP/IP , TLS , Crypto, Bluetooth, GNSS" (plus
# its own underline connector). The "ABSTRACTION LAYERS" caption that
# floated above that area was repositioned to sit above the new band.

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# --- grab style/formatting donors before we delete anything -----------
# "Straight Connector 82" (id 83) is the accent1 underline rule used
# under the abstraction-layer bands; duplicate it to build the new rule.
$connDonor = Get-ShapeById $s 83
$newConn = $connDonor.Duplicate().Item(1)

# "Rounded Rectangle 87" (id 89) is one of the small "Drivers" pills
# that already has the exact bodyPr/pPr/rPr shape we need for the new
# band's text (no stray lnSpc / inset overrides); duplicate it rather
# than the old id-95 band so we inherit that clean formatting.
$rrectDonor = Get-ShapeById $s 89
$newRRect = $rrectDonor.Duplicate().Item(1)

# "Rectangle 108" (id 109) is the "ABSTRACTION LAYERS" caption itself;
# duplicate it so the copy can be repositioned, and delete the original.
$rectDonor = Get-ShapeById $s 109
$newRect = $rectDonor.Duplicate().Item(1)

# --- delete the shapes this edit removes -------------------------------
(Get-ShapeById $s 82).Delete()   # Rounded Rectangle 73 - "Graphical Engine"
(Get-ShapeById $s 83).Delete()   # Straight Connector 82 (old underline)
(Get-ShapeById $s 95).Delete()   # Rounded Rectangle 73 - old "C stacks..." band
(Get-ShapeById $s 109).Delete()  # Rectangle 108 - old "ABSTRACTION LAYERS" position

# --- new underline connector (id 13 in the authored file) --------------
$newConn.Name = "Straight Connector 12"
$newConn.Left = 455.5165354330709
$newConn.Top = 318.87544451495796
$newConn.Width = 282.20772035013033
$newConn.Height = 0.0

# --- new merged band: "C stacks TCP/IP , TLS , Crypto, Bluetooth, GNSS" (id 15) --
$newRRect.Name = "Rounded Rectangle 27"
$newRRect.Left = 453.36457837322564
$newRRect.Top = 327.1025273144339
$newRRect.Width = 286.5071653543307
$newRRect.Height = 22.753858267716534
$newRRect.Adjustments.Item(1) = 0.13128
$newRRect.TextFrame.TextRange.Text = "C stacks TCP/IP , TLS , Crypto, Bluetooth, GNSS"

# --- "ABSTRACTION LAYERS" caption, moved above the new band (id 18) ----
$newRect.Name = "Rectangle 17"
$newRect.Left = 519.8714325816237
$newRect.Top = 314.9388303417296
$newRect.Width = 147.78315151364785
$newRect.Height = 10.203700787401575
